$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (Jengibre @ Terminal La Palmera de La Serena) was
# added to the daily feed. It belongs chronologically right above the
# existing row 97, so insert a fresh row there and push rows 97:111 down to
# 98:112 (preserving all of their data untouched).
$ws.Rows.Item(97).EntireRow.Insert()

# Populate the newly inserted row 97 with the new record's data.
$ws.Cells.Item(97, 1).Value = 8
$ws.Cells.Item(97, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(97, 3).Value = "Coquimbo"
$ws.Cells.Item(97, 4).Value = 44988
$ws.Cells.Item(97, 5).Value = 4
$ws.Cells.Item(97, 6).Value = 100114007
$ws.Cells.Item(97, 7).Value = "Jengibre"
$ws.Cells.Item(97, 8).Value = "Sin especificar"
$ws.Cells.Item(97, 9).Value = "Primera"
$ws.Cells.Item(97, 10).Value = 400
$ws.Cells.Item(97, 11).Value = 19000
$ws.Cells.Item(97, 12).Value = 20000
$ws.Cells.Item(97, 13).Value = 19500
$ws.Cells.Item(97, 14).Value = "$/caja 13 kilos"
$ws.Cells.Item(97, 15).Value = "Perú"
$ws.Cells.Item(97, 16).Value = 1500
$ws.Cells.Item(97, 17).Value = 13
$ws.Cells.Item(97, 18).Value = "Hortaliza"
